$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.0031744077509557, 1.252230448246237, 0.31518458944580485, 3.08106302114634)"
$ws.Range("C2").Value = "JSU(-1.4129619153440953, 1.1160040177193418, 2.3470392309861445, 4.137621400023746)"
$ws.Range("D2").Value = "JSU(-0.8783791409560409, 1.0369608873767402, 0.5718472633795239, 2.550528342763366)"
$ws.Range("E2").Value = "NIG(1.5810429003730713, 1.250772096908582, 3.251317475695074, 5.936769123514924)"
